# Kansas Home Defensive Actions - cleaned defensive actions data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge the header cell groups in row 1 so each column gets its own label.
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# Row 1: top-level (group) header labels.
$ws.Range("A1").Value = "Player ID"
$ws.Range("B1").Value = "Player"
$ws.Range("C1").Value = "#"
$ws.Range("D1").Value = "Nation"
$ws.Range("E1").Value = "Pos"
$ws.Range("F1").Value = "Age"
$ws.Range("G1").Value = "90s"
$ws.Range("H1").Value = "Tkl"
$ws.Range("I1").Value = "TklW"
$ws.Range("J1").Value = "Def 3rd"
$ws.Range("K1").Value = "Mid 3rd"
$ws.Range("L1").Value = "Att 3rd"
$ws.Range("M1").Value = "Cha"
$ws.Range("N1").Value = "Att"
$ws.Range("O1").Value = "Tkl%"
$ws.Range("P1").Value = "Lost"
$ws.Range("Q1").Value = "Blocks"
$ws.Range("R1").Value = "Sh"
$ws.Range("S1").Value = "Pass"
$ws.Range("T1").Value = "Int"
$ws.Range("U1").Value = "Tkl+Int"
$ws.Range("V1").Value = "Clr"
$ws.Range("W1").Value = "Err"

# Row 2 keeps its existing sub-header labels, but is now a hidden helper row.
$ws.Rows.Item(2).Hidden = $true

# New blank hidden spacer row between the headers and the data.
$ws.Rows.Item(3).Hidden = $true

# A couple of previously-omitted zero values now present in the data rows.
$ws.Range("O5").Value = 0
$ws.Range("O18").Value = 0

# The trailing summary row is now hidden too.
$ws.Rows.Item(19).Hidden = $true

# Restore the selection to match the saved workbook state.
$ws.Range("O20").Select()
